# Expense_Tracker.xlsx - "Minor changes" commit
# - Reshuffles Category/Amount pairs that were entered out of order within several
#   same-day groups (values move between sibling rows; the Date/Notes columns are untouched).
# - Inserts a new test row (Restaurant, 30, "test6") into the 2025/04/29 "test" block,
#   pushing the table from 206 to 207 data rows, and re-sorts that block into date order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at 202 (shifts old rows 202-206 down to 203-207).
$ws.Rows.Item(202).Insert()

# --- Fix up the Category/Amount pairs that were out of order within their date group ---
$ws.Cells.Item(14, 2).Value = "Barber"
$ws.Cells.Item(14, 3).Value = 89.78
$ws.Cells.Item(15, 2).Value = "Transportation"
$ws.Cells.Item(15, 3).Value = 119.26
$ws.Cells.Item(16, 2).Value = "Phone"
$ws.Cells.Item(16, 3).Value = 105.95
$ws.Cells.Item(27, 2).Value = "Laundry"
$ws.Cells.Item(27, 3).Value = 76.5
$ws.Cells.Item(28, 2).Value = "Phone"
$ws.Cells.Item(28, 3).Value = 129.89
$ws.Cells.Item(52, 2).Value = "Shopping"
$ws.Cells.Item(52, 3).Value = 177.51
$ws.Cells.Item(53, 2).Value = "Shopping"
$ws.Cells.Item(53, 3).Value = 169.37
$ws.Cells.Item(54, 2).Value = "Restaurant"
$ws.Cells.Item(54, 3).Value = 42.06
$ws.Cells.Item(55, 2).Value = "Entertainment"
$ws.Cells.Item(55, 3).Value = 91.19
$ws.Cells.Item(77, 2).Value = "Groceries"
$ws.Cells.Item(77, 3).Value = 36.96
$ws.Cells.Item(78, 2).Value = "Restaurant"
$ws.Cells.Item(78, 3).Value = 25.39
$ws.Cells.Item(79, 2).Value = "Snacks"
$ws.Cells.Item(79, 3).Value = 19.68
$ws.Cells.Item(80, 2).Value = "Restaurant"
$ws.Cells.Item(80, 3).Value = 49.5
$ws.Cells.Item(81, 2).Value = "Toters"
$ws.Cells.Item(81, 3).Value = 40.23
$ws.Cells.Item(88, 2).Value = "Shopping"
$ws.Cells.Item(88, 3).Value = 119.13
$ws.Cells.Item(89, 2).Value = "Shopping"
$ws.Cells.Item(89, 3).Value = 17.42
$ws.Cells.Item(90, 2).Value = "Transportation"
$ws.Cells.Item(90, 3).Value = 71.19
$ws.Cells.Item(91, 2).Value = "Shopping"
$ws.Cells.Item(91, 3).Value = 50.47
$ws.Cells.Item(98, 2).Value = "Laundry"
$ws.Cells.Item(98, 3).Value = 182.12
$ws.Cells.Item(99, 2).Value = "Restaurant"
$ws.Cells.Item(99, 3).Value = 28.46
$ws.Cells.Item(100, 2).Value = "Groceries"
$ws.Cells.Item(100, 3).Value = 80.29000000000001
$ws.Cells.Item(101, 2).Value = "Toters"
$ws.Cells.Item(101, 3).Value = 52.66
$ws.Cells.Item(102, 2).Value = "Snacks"
$ws.Cells.Item(102, 3).Value = 12.83
$ws.Cells.Item(129, 2).Value = "Phone"
$ws.Cells.Item(129, 3).Value = 75.84999999999999
$ws.Cells.Item(130, 2).Value = "Laundry"
$ws.Cells.Item(130, 3).Value = 151.66
$ws.Cells.Item(131, 2).Value = "Transportation"
$ws.Cells.Item(131, 3).Value = 80.59
$ws.Cells.Item(132, 2).Value = "Transportation"
$ws.Cells.Item(132, 3).Value = 62.92
$ws.Cells.Item(140, 2).Value = "Restaurant"
$ws.Cells.Item(140, 3).Value = 83.48
$ws.Cells.Item(141, 2).Value = "Toters"
$ws.Cells.Item(141, 3).Value = 169.93
$ws.Cells.Item(142, 2).Value = "Groceries"
$ws.Cells.Item(142, 3).Value = 21.07
$ws.Cells.Item(143, 2).Value = "Laundry"
$ws.Cells.Item(143, 3).Value = 62.67
$ws.Cells.Item(144, 2).Value = "Restaurant"
$ws.Cells.Item(144, 3).Value = 45.66
$ws.Cells.Item(156, 2).Value = "Laundry"
$ws.Cells.Item(156, 3).Value = 121.27
$ws.Cells.Item(157, 2).Value = "Phone"
$ws.Cells.Item(157, 3).Value = 34.21
$ws.Cells.Item(158, 2).Value = "Shopping"
$ws.Cells.Item(158, 3).Value = 180.79
$ws.Cells.Item(168, 2).Value = "Entertainment"
$ws.Cells.Item(168, 3).Value = 38.46
$ws.Cells.Item(169, 2).Value = "Barber"
$ws.Cells.Item(169, 3).Value = 12.48
$ws.Cells.Item(170, 2).Value = "Toters"
$ws.Cells.Item(170, 3).Value = 171.59
$ws.Cells.Item(171, 2).Value = "Shopping"
$ws.Cells.Item(171, 3).Value = 106.97
$ws.Cells.Item(180, 2).Value = "Transportation"
$ws.Cells.Item(180, 3).Value = 69.01000000000001
$ws.Cells.Item(181, 2).Value = "Barber"
$ws.Cells.Item(181, 3).Value = 189.01
$ws.Cells.Item(182, 2).Value = "Restaurant"
$ws.Cells.Item(182, 3).Value = 14.78
$ws.Cells.Item(188, 2).Value = "Restaurant"
$ws.Cells.Item(188, 3).Value = 183.49
$ws.Cells.Item(189, 2).Value = "Entertainment"
$ws.Cells.Item(189, 3).Value = 199.43
$ws.Cells.Item(190, 2).Value = "Groceries"
$ws.Cells.Item(190, 3).Value = 171.79
$ws.Cells.Item(191, 2).Value = "Transportation"
$ws.Cells.Item(191, 3).Value = 69.18000000000001
$ws.Cells.Item(192, 2).Value = "Toters"
$ws.Cells.Item(192, 3).Value = 92.91

# --- Rebuild the 2025/04/29 "test" rows (200-207) in their corrected order ---
$ws.Cells.Item(200, 2).Value = "Restaurant"
$ws.Cells.Item(200, 3).Value = 20
$ws.Cells.Item(200, 4).Value = "test5"
$ws.Cells.Item(201, 4).Value = "test3"
$ws.Cells.Item(202, 1).Value = "'2025/04/29"
$ws.Cells.Item(202, 2).Value = "Restaurant"
$ws.Cells.Item(202, 3).Value = 30
$ws.Cells.Item(202, 4).Value = "test6"
$ws.Cells.Item(203, 4).Value = "test4"
$ws.Cells.Item(204, 2).Value = "Barber"
$ws.Cells.Item(204, 3).Value = 69.84999999999999
$ws.Cells.Item(204, 4).Value = ""
$ws.Cells.Item(205, 4).Value = "test"
